# Update FFXIV leve-crafting profit figures (columns H-N) per refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1025.4706
$ws.Range("I19").Value = 649.1818
$ws.Range("K19").Value = 649.1818
$ws.Range("M19").Value = -474.1818
# Row 99
$ws.Range("H99").Value = 204.66667
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = ""
# Row 132
$ws.Range("H132").Value = 7766.7896
$ws.Range("I132").Value = 7723.75
$ws.Range("K132").Value = 23171.25
$ws.Range("M132").Value = -20641.25
# Row 135
$ws.Range("H135").Value = 1214.7826
$ws.Range("I135").Value = 1179.091
$ws.Range("K135").Value = 10611.819
$ws.Range("M135").Value = -8076.819
# Row 137
$ws.Range("H137").Value = 1750.5366
$ws.Range("J137").Value = 1885.4286
$ws.Range("L137").Value = 5656.2858
$ws.Range("N137").Value = -10756.2858
# Row 138
$ws.Range("H138").Value = 5939.9443
$ws.Range("I138").Value = 5108.5386
$ws.Range("J138").Value = 6409.8696
$ws.Range("K138").Value = 15325.6158
$ws.Range("L138").Value = 19229.6088
$ws.Range("M138").Value = -10185.6158
$ws.Range("N138").Value = -29509.6088

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 1571.4286
$ws.Range("J4").Value = 1974.75
$ws.Range("L4").Value = 1974.75
$ws.Range("N4").Value = -2206.75
# Row 5
$ws.Range("H5").Value = 658.3333
$ws.Range("I5").Value = 683.3333
$ws.Range("J5").Value = 633.3333
$ws.Range("K5").Value = 683.3333
$ws.Range("L5").Value = 633.3333
$ws.Range("M5").Value = -571.3333
$ws.Range("N5").Value = -857.3333
# Row 61
$ws.Range("H61").Value = 3973217.2
$ws.Range("I61").Value = 5212483
$ws.Range("J61").Value = 7566.2
$ws.Range("K61").Value = 5212483
$ws.Range("L61").Value = 7566.2
$ws.Range("M61").Value = -5212271
$ws.Range("N61").Value = -7990.2
# Row 122
$ws.Range("H122").Value = 2251.5
$ws.Range("I122").Value = 1678.9375
$ws.Range("J122").Value = 3778.3333
$ws.Range("K122").Value = 5036.8125
$ws.Range("L122").Value = 11334.9999
$ws.Range("M122").Value = -2586.8125
$ws.Range("N122").Value = -16234.9999
# Row 132
$ws.Range("H132").Value = 4125.6514
$ws.Range("I132").Value = 3525.9722
$ws.Range("K132").Value = 10577.9166
$ws.Range("M132").Value = -8047.9166
# Row 134
$ws.Range("H134").Value = 149666.33
$ws.Range("J134").Value = 149666.33
$ws.Range("L134").Value = 149666.33
$ws.Range("N134").Value = -159806.33
# Row 136
$ws.Range("H136").Value = 3973217.2
$ws.Range("I136").Value = 5212483
$ws.Range("J136").Value = 7566.2
$ws.Range("K136").Value = 15637449
$ws.Range("L136").Value = 22698.6
$ws.Range("M136").Value = -15634899
$ws.Range("N136").Value = -27798.6

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 658.3333
$ws.Range("I4").Value = 683.3333
$ws.Range("J4").Value = 633.3333
$ws.Range("K4").Value = 683.3333
$ws.Range("L4").Value = 633.3333
$ws.Range("M4").Value = -568.3333
$ws.Range("N4").Value = -863.3333
# Row 64
$ws.Range("H64").Value = 729.4
$ws.Range("J64").Value = 602.7143
$ws.Range("L64").Value = 602.7143
$ws.Range("N64").Value = -1052.7143
# Row 67
$ws.Range("H67").Value = 729.4
$ws.Range("J67").Value = 602.7143
$ws.Range("L67").Value = 602.7143
$ws.Range("N67").Value = -2162.7143
# Row 134
$ws.Range("H134").Value = 4829.45
$ws.Range("I134").Value = 5747.3076
$ws.Range("K134").Value = 17241.9228
$ws.Range("M134").Value = -14706.9228

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6855.722
$ws.Range("I31").Value = 5322.6665
$ws.Range("K31").Value = 5322.6665
$ws.Range("M31").Value = -5027.6665
# Row 34
$ws.Range("H34").Value = 6855.722
$ws.Range("I34").Value = 5322.6665
$ws.Range("K34").Value = 5322.6665
$ws.Range("M34").Value = -5120.6665
# Row 59
$ws.Range("H59").Value = 28386
$ws.Range("J59").Value = 29232.5
$ws.Range("L59").Value = 29232.5
$ws.Range("N59").Value = -31522.5
# Row 86
$ws.Range("H86").Value = 4417.4736
$ws.Range("I86").Value = 3720.75
$ws.Range("J86").Value = 8133.3335
$ws.Range("K86").Value = 3720.75
$ws.Range("L86").Value = 8133.3335
$ws.Range("M86").Value = -2597.75
$ws.Range("N86").Value = -10379.3335
# Row 89
$ws.Range("H89").Value = 4417.4736
$ws.Range("I89").Value = 3720.75
$ws.Range("J89").Value = 8133.3335
$ws.Range("K89").Value = 18603.75
$ws.Range("L89").Value = 40666.6675
$ws.Range("M89").Value = -12987.75
$ws.Range("N89").Value = -51898.6675
# Row 94
$ws.Range("H94").Value = 2482.4285
$ws.Range("I94").Value = 2144.5
$ws.Range("J94").Value = 2617.6
$ws.Range("K94").Value = 2144.5
$ws.Range("L94").Value = 2617.6
$ws.Range("M94").Value = -1693.5
$ws.Range("N94").Value = -3519.6
# Row 124
$ws.Range("H124").Value = 47663
$ws.Range("J124").Value = 47663
$ws.Range("L124").Value = 47663
$ws.Range("N124").Value = -52573
# Row 134
$ws.Range("H134").Value = 6741.5713
$ws.Range("I134").Value = 4235.875
$ws.Range("J134").Value = 10082.5
$ws.Range("K134").Value = 12707.625
$ws.Range("L134").Value = 30247.5
$ws.Range("M134").Value = -10172.625
$ws.Range("N134").Value = -35317.5

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 433.33334
$ws.Range("J34").Value = 500
$ws.Range("L34").Value = 1500
$ws.Range("N34").Value = -1668
# Row 39
$ws.Range("H39").Value = 1000
$ws.Range("J39").Value = 1000
$ws.Range("L39").Value = 3000
$ws.Range("N39").Value = -3588
# Row 54
$ws.Range("H54").Value = 12596
$ws.Range("J54").Value = 12596
$ws.Range("L54").Value = 37788
$ws.Range("N54").Value = -38906
# Row 55
$ws.Range("H55").Value = 3999
$ws.Range("J55").Value = 3999
$ws.Range("L55").Value = 11997
$ws.Range("N55").Value = -12351
# Row 122
$ws.Range("H122").Value = 373.31818
$ws.Range("J122").Value = 359.22223
$ws.Range("L122").Value = 3233.00007
$ws.Range("N122").Value = -8133.00007
# Row 129
$ws.Range("H129").Value = 556735.7
$ws.Range("J129").Value = 722554.2
$ws.Range("L129").Value = 2167662.6
$ws.Range("N129").Value = -2177662.6
# Row 131
$ws.Range("H131").Value = 35715136
$ws.Range("I131").Value = 50000396
$ws.Range("K131").Value = 150001188
$ws.Range("M131").Value = -149996148

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1610.3846
$ws.Range("I132").Value = 1610.3846
$ws.Range("K132").Value = 4831.1538
$ws.Range("M132").Value = -2301.1538

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 8584.388999999999
$ws.Range("I132").Value = 9617.5
$ws.Range("J132").Value = 4968.5
$ws.Range("K132").Value = 28852.5
$ws.Range("L132").Value = 14905.5
$ws.Range("M132").Value = -26322.5
$ws.Range("N132").Value = -19965.5
# Row 136
$ws.Range("H136").Value = 4413.316
$ws.Range("I136").Value = 3390.6
$ws.Range("J136").Value = 8248.5
$ws.Range("K136").Value = 10171.8
$ws.Range("L136").Value = 24745.5
$ws.Range("M136").Value = -7621.799999999999
$ws.Range("N136").Value = -29845.5
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""
# Row 140
$ws.Range("H140").Value = 95980
$ws.Range("J140").Value = 95980
$ws.Range("L140").Value = 95980
$ws.Range("N140").Value = -106340

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 4488.298
$ws.Range("I132").Value = 3931.325
$ws.Range("J132").Value = 7671
$ws.Range("K132").Value = 11793.975
$ws.Range("L132").Value = 23013
$ws.Range("M132").Value = -9263.974999999999
$ws.Range("N132").Value = -28073
# Row 136
$ws.Range("H136").Value = 4417.116
$ws.Range("J136").Value = 8211.286
$ws.Range("L136").Value = 24633.858
$ws.Range("N136").Value = -29733.858

